# Fruta / hortaliza, semanal
# Insert two new weekly price rows at the top of the data table (rows 2-3),
# pushing all existing records down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 2 (the first data row).
$ws.Rows("2:3").Insert()

# The D column (Fecha) uses a date/time number format; the freshly inserted
# rows don't inherit it automatically, so set it explicitly before writing
# the date values.
$ws.Range("D2:D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Row 2: new "Primera" quality record dated 2023-11-07 ---
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2").Value = "Ñuble"
$ws.Range("D2").Value = 45237
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100107
$ws.Range("H2").Value = "Otros"
$ws.Range("I2").Value = 100107002
$ws.Range("J2").Value = "Chirimoya"
$ws.Range("K2").Value = "Cultivar IV Región"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 22000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 22000
$ws.Range("Q2").Value = "`$/bandeja 10 kilos"
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 2200
$ws.Range("T2").Value = 10

# --- Row 3: new "Segunda" quality record, same date ---
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 45237
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100107
$ws.Range("H3").Value = "Otros"
$ws.Range("I3").Value = 100107002
$ws.Range("J3").Value = "Chirimoya"
$ws.Range("K3").Value = "Cultivar IV Región"
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("Q3").Value = "`$/bandeja 10 kilos"
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 1800
$ws.Range("T3").Value = 10
